# Generate Report for Handback
# The 8722a735-fc9e-4e87-ac39-ac760696e3a7.md file has now been handed back
# in sync with en-US for both zh-cn and de-de locales. Update the status,
# the Latest Handback DateTime, and clear the stale Error Detail message
# on each locale sheet, and reflect the new status on the Overview sheet.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-10-17 14:37:19"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-10-17 14:38:00"
$wsDeDe.Range("P3").Value = ""
